$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values - force text format to preserve exact string representation
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.725.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.800.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.441.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.787.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.801.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.693"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.947.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.741.35"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.1000"
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "395.01"
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  +1.39%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").Value = "  +2.83%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  +8.62%  "
$ws.Range("E50").Value = "  +7.54%  "
$ws.Range("E51").Value = "  +0.57%  "
